$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 3 -> HK UMP45 200mm .45ACP Flanged barrel
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "hk_ump45_200mm_flanged_barrel"
$ws.Range("B3").Value = "HK UMP45 200mm .45ACP Flanged"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.26
$ws.Range("E3").Value = -2
$ws.Range("F3").Value = -1
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = -15
$ws.Range("M3").Value = 0
$ws.Range("Q3").Value = 7.8740199999999998

# ---------------------------------------------------------------------------
# 2. Row 4 -> HDPS UMP45 200mm .45ACP Threaded barrel
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "hdps_ump45_200mm_threaded_barrel"
$ws.Range("B4").Value = "HDPS UMP45 200mm .45ACP Threaded"
$ws.Range("C4").Value = -2
$ws.Range("D4").Value = 0.27
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("H4").Value = 0.08
$ws.Range("I4").ClearContents()
$ws.Range("J4").Value = -21
$ws.Range("M4").Value = 750
$ws.Range("Q4").Value = 7.8740199999999998

# ---------------------------------------------------------------------------
# 3. Row 5 -> OMEGA UMP45 200mm .45ACP Tri-Lug barrel
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "omega_ump45_200mm_trilug_barrel"
$ws.Range("B5").Value = "OMEGA UMP45 200mm .45ACP Tri-Lug"
$ws.Range("C5").Value = -1
$ws.Range("D5").Value = 0.26
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = -2
$ws.Range("H5").Value = 0.05
$ws.Range("I5").ClearContents()
$ws.Range("J5").Value = -18
$ws.Range("M5").Value = 300
$ws.Range("Q5").Value = 7.8740199999999998

# ---------------------------------------------------------------------------
# 4. Row 7 -> Kriss Vector .45ACP 170mm barrel (updated stats)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "vector_45acp_170mm_barrel"
$ws.Range("B7").Value = "Kriss Vector .45ACP 170mm"
$ws.Range("C7").Value = -1
$ws.Range("D7").Value = 0.23
$ws.Range("E7").Value = -1
$ws.Range("F7").Value = -1
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = -0.05
$ws.Range("I7").Value = 0.03
$ws.Range("J7").Value = -72
$ws.Range("K7").Value = -0.1
$ws.Range("L7").ClearContents()
$ws.Range("M7").Value = 750

# ---------------------------------------------------------------------------
# 5. Row 8 -> Kriss Vector .45ACP 140mm barrel (updated stats)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "vector_45acp_140mm_barrel"
$ws.Range("B8").Value = "Kriss Vector .45ACP 140mm"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0.2
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = -99
$ws.Range("K8").Value = 0
$ws.Range("L8").ClearContents()
$ws.Range("M8").Value = 800

# ---------------------------------------------------------------------------
# 6. Strip the "applied fill" style left over from the old layout so the
#    cellXfs table collapses back down, then drop the now-empty rows
#    9:38 tail down to row 22 and clear out the leftover filler cells in
#    rows 9:22 (only the helper formula in column N should remain).
# ---------------------------------------------------------------------------
$ws.Range("A1:N2").ClearFormats()
$ws.Range("N3:N6").ClearFormats()
$ws.Range("A7:N8").ClearFormats()
$ws.Range("A9:N22").ClearFormats()
$ws.Range("A9:M22").ClearContents()

$ws.Rows("23:38").Delete()

$ws.Range("A1:B1").ClearContents()
$ws.Range("D1:N1").ClearContents()

# ---------------------------------------------------------------------------
# 7. Column widths: drop the old column-B-only width and set new widths for
#    column A and column B.
# ---------------------------------------------------------------------------
$ws.Columns("B").ClearFormats()
$ws.Columns("A").ColumnWidth = 25.140625
$ws.Columns("B").ColumnWidth = 38.28515625

# ---------------------------------------------------------------------------
# 8. Selection moves to F7.
# ---------------------------------------------------------------------------
$ws.Range("F7").Select()
